# Scheduled runner refresh: update market-board derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on each Leve sheet with
# freshly pulled values.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 498.2353
$ws.Range("I6").Value = 111.333336
$ws.Range("J6").Value = 3400
$ws.Range("K6").Value = 334.000008
$ws.Range("L6").Value = 10200
$ws.Range("M6").Value = -222.000008
$ws.Range("N6").Value = -10424
$ws.Range("H8").Value = 876.55554
$ws.Range("I8").Value = 41.285713
$ws.Range("K8").Value = 123.857139
$ws.Range("M8").Value = 15.142861
$ws.Range("H38").Value = 3531.4285
$ws.Range("I38").Value = 92.5
$ws.Range("K38").Value = 277.5
$ws.Range("M38").Value = 94.5
$ws.Range("H132").Value = 76711.016
$ws.Range("I132").Value = 85927.74000000001
$ws.Range("J132").Value = 6049.5
$ws.Range("K132").Value = 257783.22
$ws.Range("L132").Value = 18148.5
$ws.Range("M132").Value = -255253.22
$ws.Range("N132").Value = -23208.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 10190.8
$ws.Range("I6").Value = 8666.666999999999
$ws.Range("J6").Value = 12477
$ws.Range("K6").Value = 8666.666999999999
$ws.Range("L6").Value = 12477
$ws.Range("M6").Value = -8493.666999999999
$ws.Range("N6").Value = -12823
$ws.Range("H36").Value = 30000
$ws.Range("J36").Value = 30000
$ws.Range("L36").Value = 30000
$ws.Range("N36").Value = -30692
$ws.Range("H45").Value = 980.86957
$ws.Range("I45").Value = 769.0909
$ws.Range("J45").Value = 1175
$ws.Range("K45").Value = 769.0909
$ws.Range("L45").Value = 1175
$ws.Range("M45").Value = -392.0909
$ws.Range("N45").Value = -1929
$ws.Range("H139").Value = 45776.875
$ws.Range("J139").Value = 45776.875
$ws.Range("L139").Value = 45776.875
$ws.Range("N139").Value = -56056.875

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 41045
$ws.Range("J138").Value = 41045
$ws.Range("L138").Value = 41045
$ws.Range("N138").Value = -51325

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2720
$ws.Range("I86").Value = 2626.75
$ws.Range("J86").Value = 2968.6667
$ws.Range("K86").Value = 2626.75
$ws.Range("L86").Value = 2968.6667
$ws.Range("M86").Value = -1503.75
$ws.Range("N86").Value = -5214.6667
$ws.Range("H89").Value = 2720
$ws.Range("I89").Value = 2626.75
$ws.Range("J89").Value = 2968.6667
$ws.Range("K89").Value = 13133.75
$ws.Range("L89").Value = 14843.3335
$ws.Range("M89").Value = -7517.75
$ws.Range("N89").Value = -26075.3335
$ws.Range("H132").Value = 2397.8096
$ws.Range("I132").Value = 977.3913
$ws.Range("K132").Value = 2932.1739
$ws.Range("M132").Value = -402.1738999999998
$ws.Range("H140").Value = 39190
$ws.Range("J140").Value = 39190
$ws.Range("L140").Value = 39190
$ws.Range("N140").Value = -49550
$ws.Range("H141").Value = 24378.262
$ws.Range("J141").Value = 24378.262
$ws.Range("L141").Value = 24378.262
$ws.Range("N141").Value = -34738.262

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1879787.4
$ws.Range("I2").Value = 70.25
$ws.Range("J2").Value = 2381045.2
$ws.Range("K2").Value = 421.5
$ws.Range("L2").Value = 14286271.2
$ws.Range("M2").Value = -308.5
$ws.Range("N2").Value = -14286497.2
$ws.Range("H117").Value = 12200
$ws.Range("I117").Value = 600
$ws.Range("J117").Value = 16066.667
$ws.Range("K117").Value = 1800
$ws.Range("L117").Value = 48200.001
$ws.Range("M117").Value = 1642
$ws.Range("N117").Value = -55084.001
$ws.Range("H131").Value = 6411146
$ws.Range("I131").Value = 100002210
$ws.Range("J131").Value = 799.35614
$ws.Range("K131").Value = 300006630
$ws.Range("L131").Value = 2398.06842
$ws.Range("M131").Value = -300001590
$ws.Range("N131").Value = -12478.06842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 15000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 15000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 15000
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -15226
$ws.Range("H16").Value = 15000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 15000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 15000
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -15500
$ws.Range("H70").Value = 6433.6875
$ws.Range("I70").Value = 5964.871
$ws.Range("J70").Value = 7288.5884
$ws.Range("K70").Value = 5964.871
$ws.Range("L70").Value = 7288.5884
$ws.Range("M70").Value = -5694.871
$ws.Range("N70").Value = -7828.5884
$ws.Range("H73").Value = 6433.6875
$ws.Range("I73").Value = 5964.871
$ws.Range("J73").Value = 7288.5884
$ws.Range("K73").Value = 5964.871
$ws.Range("L73").Value = 7288.5884
$ws.Range("M73").Value = -5028.871
$ws.Range("N73").Value = -9160.588400000001
$ws.Range("H122").Value = 3864.9
$ws.Range("I122").Value = 2354.889
$ws.Range("J122").Value = 5100.364
$ws.Range("K122").Value = 7064.667
$ws.Range("L122").Value = 15301.092
$ws.Range("M122").Value = -4614.667
$ws.Range("N122").Value = -20201.092
$ws.Range("H132").Value = 4326.294
$ws.Range("I132").Value = 1499
$ws.Range("J132").Value = 4503
$ws.Range("K132").Value = 4497
$ws.Range("L132").Value = 13509
$ws.Range("M132").Value = -1967
$ws.Range("N132").Value = -18569
$ws.Range("H140").Value = 42794.547
$ws.Range("J140").Value = 42794.547
$ws.Range("L140").Value = 42794.547
$ws.Range("N140").Value = -53154.547

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3123.2307
$ws.Range("I22").Value = 2333.3333
$ws.Range("J22").Value = 3360.2
$ws.Range("K22").Value = 2333.3333
$ws.Range("L22").Value = 3360.2
$ws.Range("M22").Value = -2038.3333
$ws.Range("N22").Value = -3950.2
$ws.Range("H27").Value = 3123.2307
$ws.Range("I27").Value = 2333.3333
$ws.Range("J27").Value = 3360.2
$ws.Range("K27").Value = 2333.3333
$ws.Range("L27").Value = 3360.2
$ws.Range("M27").Value = -2226.3333
$ws.Range("N27").Value = -3574.2
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H46").Value = 1752.04
$ws.Range("I46").Value = 1470.05
$ws.Range("K46").Value = 1470.05
$ws.Range("M46").Value = -1282.05
$ws.Range("H139").Value = 39773
$ws.Range("J139").Value = 39773
$ws.Range("L139").Value = 39773
$ws.Range("N139").Value = -50053
$ws.Range("H140").Value = 69089.7
$ws.Range("J140").Value = 69089.7
$ws.Range("L140").Value = 69089.7
$ws.Range("N140").Value = -79449.7
$ws.Range("H141").Value = 40922.633
$ws.Range("J141").Value = 40922.633
$ws.Range("L141").Value = 40922.633
$ws.Range("N141").Value = -51282.633

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 29000
$ws.Range("J31").Value = 29000
$ws.Range("L31").Value = 29000
$ws.Range("N31").Value = -29696
$ws.Range("H107").Value = 597.381
$ws.Range("I107").Value = 491.44446
$ws.Range("J107").Value = 1233
$ws.Range("K107").Value = 1474.33338
$ws.Range("L107").Value = 3699
$ws.Range("M107").Value = 445.66662
$ws.Range("N107").Value = -7539
$ws.Range("H122").Value = 3895.75
$ws.Range("I122").Value = 1992.9166
$ws.Range("K122").Value = 5978.7498
$ws.Range("M122").Value = -3528.7498
$ws.Range("H126").Value = 3650.0588
$ws.Range("I126").Value = 3013.6667
$ws.Range("J126").Value = 3997.182
$ws.Range("K126").Value = 9041.000100000001
$ws.Range("L126").Value = 11991.546
$ws.Range("M126").Value = -6571.000100000001
$ws.Range("N126").Value = -16931.546
$ws.Range("H138").Value = 52682.5
$ws.Range("J138").Value = 52682.5
$ws.Range("L138").Value = 52682.5
$ws.Range("N138").Value = -62962.5
$ws.Range("H139").Value = 41002.918
$ws.Range("J139").Value = 41018.26
$ws.Range("L139").Value = 41018.26
$ws.Range("N139").Value = -51298.26
$ws.Range("H141").Value = 43877.918
$ws.Range("J141").Value = 43877.918
$ws.Range("L141").Value = 43877.918
$ws.Range("N141").Value = -54237.918
